$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $cell = $ws.Range($addr)
    # Coinranking price strings sometimes look like plain numbers
    # (e.g. "1.00", "0.140", "49.70"). Force text format first so
    # Excel keeps the literal digits instead of normalising them
    # into a number (which would drop trailing/duplicate zeros).
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

$ws.Range('D2').Value = '47.212.94'
$ws.Range('E2').Value = '  -0.56%  '

$ws.Range('D3').Value = '2.484.68'
$ws.Range('E3').Value = '  -0.41%  '

$ws.Range('E4').Value = '  +0.07%  '

Set-TextCell 'D5' '321.03'
$ws.Range('E5').Value = '  -0.71%  '

Set-TextCell 'D6' '108.06'
$ws.Range('E6').Value = '  +2.66%  '

Set-TextCell 'D7' '0.522'
$ws.Range('E7').Value = '  +0.16%  '

Set-TextCell 'D8' '1.00'
$ws.Range('E8').Value = '  +0.05%  '

$ws.Range('E9').Value = '  -0.07%  '

Set-TextCell 'D10' '39.06'
$ws.Range('E10').Value = '  +4.03%  '

Set-TextCell 'D11' '0.0809'
$ws.Range('E11').Value = '  -0.34%  '

Set-TextCell 'D13' '18.34'
$ws.Range('E13').Value = '  +0.06%  '

Set-TextCell 'D14' '7.17'
$ws.Range('E14').Value = '  -0.34%  '

$ws.Range('D15').Value = '2.873.95'
$ws.Range('E15').Value = '  -0.20%  '

$ws.Range('D16').Value = '2.485.38'
$ws.Range('E16').Value = '  -0.55%  '

Set-TextCell 'D17' '0.842'
$ws.Range('E17').Value = '  -0.03%  '

$ws.Range('D18').Value = '47.108.88'
$ws.Range('E18').Value = '  -0.48%  '

Set-TextCell 'D19' '13.31'
$ws.Range('E19').Value = '  +4.46%  '

Set-TextCell 'D20' '6.61'
$ws.Range('E20').Value = '  +0.95%  '

$ws.Range('E21').Value = '  +0.87%  '

Set-TextCell 'D22' '2.74'
$ws.Range('E22').Value = '  +14.55%  '

Set-TextCell 'D23' '70.51'
$ws.Range('E23').Value = '  -0.27%  '

Set-TextCell 'D24' '245.76'
$ws.Range('E24').Value = '  -1.98%  '

Set-TextCell 'D25' '2.56'
$ws.Range('E25').Value = '  -0.34%  '

$ws.Range('E26').Value = '  -0.03%  '

Set-TextCell 'D27' '25.67'
$ws.Range('E27').Value = '  -1.84%  '

$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell 'D28' '9.96'
$ws.Range('E28').Value = '  -1.32%  '

$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell 'D29' '2.19'
$ws.Range('E29').Value = '  -0.52%  '

Set-TextCell 'D30' '0.140'
$ws.Range('E30').Value = '  +4.86%  '

Set-TextCell 'D31' '34.53'
$ws.Range('E31').Value = '  -1.79%  '

Set-TextCell 'D32' '49.70'
$ws.Range('E32').Value = '  +0.47%  '

Set-TextCell 'D33' '20.23'
$ws.Range('E33').Value = '  +1.10%  '

Set-TextCell 'D34' '5.32'
$ws.Range('E34').Value = '  -1.02%  '

Set-TextCell 'D35' '0.0781'
$ws.Range('E35').Value = '  -0.17%  '

$ws.Range('E36').Value = '  +0.15%  '

$ws.Range('E37').Value = '  +2.15%  '

Set-TextCell 'D38' '1.95'
$ws.Range('E38').Value = '  +0.63%  '

$ws.Range('E39').Value = '  -1.64%  '

Set-TextCell 'D40' '23.06'
$ws.Range('E40').Value = '  +7.76%  '

Set-TextCell 'D41' '0.111'
$ws.Range('E41').Value = '  +0.23%  '

Set-TextCell 'D42' '2.21'
$ws.Range('E42').Value = '  -0.94%  '

Set-TextCell 'D43' '118.71'
$ws.Range('E43').Value = '  -1.65%  '

$ws.Range('E44').Value = '  +0.12%  '

$ws.Range('D45').Value = '1.994.82'
$ws.Range('E45').Value = '  +1.71%  '

Set-TextCell 'D46' '3.01'
$ws.Range('E46').Value = '  +1.07%  '

Set-TextCell 'D47' '2.00'
$ws.Range('E47').Value = '  -4.13%  '

$ws.Range('E48').Value = '  -0.53%  '

Set-TextCell 'D49' '9.12'
$ws.Range('E49').Value = '  -1.21%  '

Set-TextCell 'D50' '5.08'
$ws.Range('E50').Value = '  -4.92%  '

Set-TextCell 'D51' '56.50'
$ws.Range('E51').Value = '  +2.85%  '

